$d = $word.ActiveDocument

# Merge the multiple split runs in the Title, Author and Abstract paragraphs
# into a single run each, by replacing the (identical) visible text via
# Find/Replace. Word's replace collapses the matched range into one run,
# removing the redundant run-splitting the authoring tool left behind.

$d.Content.Find.Execute("Answers: Laws of indices", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "Answers: Laws of indices", 2)

$d.Content.Find.Execute("Isabella Lewis, Akshat Srivastava", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "Isabella Lewis, Akshat Srivastava", 2)

$d.Content.Find.Execute("Answers to questions relating to using laws of indices.", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "Answers to questions relating to using laws of indices.", 2)
